$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 = "Marking": points awarded per right answer / per wrong answer.
# concise_ms CSV pattern update: +1 point per right answer, harsher wrong-answer penalty.
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 = "Total": recomputed from the row 10 (Right/Wrong counts) x row 11 (marking scheme).
# Right=25, Wrong=0 (see row 10) -> 25*5 + 0*-1.2
$ws.Range("B12").Value = 125
$ws.Range("C12").Value = 0 * -1.2
$ws.Range("E12").Value = "125.0/140"
